$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stats")

# Step 1: delete the old "Week 12" column (column G); remaining columns shift left
$ws.Columns.Item(7).Delete()

# Step 2: update Wins (B), PF (C) and GP (D) for each team row
$ws.Range("B2").Value = 11
$ws.Range("C2").Formula = "=1713.2+167.56"
$ws.Range("D2").Value = 12
$ws.Range("B3").Value = 8
$ws.Range("C3").Formula = "=1431.14+134.92"
$ws.Range("D3").Value = 12
$ws.Range("B4").Value = 6
$ws.Range("C4").Formula = "=1342.96+155.92"
$ws.Range("D4").Value = 12
$ws.Range("B5").Value = 7
$ws.Range("C5").Formula = "=1392.8+101.86"
$ws.Range("D5").Value = 12
$ws.Range("B6").Value = 7
$ws.Range("C6").Formula = "=1290.82+157.76"
$ws.Range("D6").Value = 12
$ws.Range("B7").Value = 8
$ws.Range("C7").Formula = "=1501.62+135.76"
$ws.Range("D7").Value = 12
$ws.Range("B8").Value = 5
$ws.Range("C8").Formula = "=1453.42+100.48"
$ws.Range("D8").Value = 12
$ws.Range("B9").Value = 5
$ws.Range("C9").Formula = "=1338.24+111.78"
$ws.Range("D9").Value = 12
$ws.Range("B10").Value = 4
$ws.Range("C10").Formula = "=1496.84+116.52"
$ws.Range("D10").Value = 12
$ws.Range("B11").Value = 5
$ws.Range("C11").Formula = "=1413.86+139.64"
$ws.Range("D11").Value = 12
$ws.Range("B12").Value = 3
$ws.Range("C12").Formula = "=1314.5+114.22"
$ws.Range("D12").Value = 12
$ws.Range("B13").Value = 3
$ws.Range("C13").Formula = "=1237.44+84.24"
$ws.Range("D13").Value = 12

# Step 3: update weekly / rest-of-season formulas (G:K) for each team row
$ws.Range("G2").Formula = "=19.39+15.8+24.14+22.36+20.58+10.95+19.12+13.51+6+8.83"
$ws.Range("H2").Formula = "=19.4+22.8+8+23.44+22.39+10.71+8.36+13.62+5.04+8.68"
$ws.Range("I2").Formula = "=19.34+16.68+28.14+12.77+24.25+10.13+18.16+20.33+6.59+8.77"
$ws.Range("J2").Formula = "=18.35+24.4+21.96+13.54+22.78+9.78+19.33+20.22+6.94+8.31"
$ws.Range("K2").Formula = "=22.23+20.07+27.7+15.19+23.51+11.58+16.7+20.44+8.77+5"
$ws.Range("G3").Formula = "=20.78+17.95+15.12+19.04+14.37+13.28+11.25+7+4.68+7.56"
$ws.Range("H3").Formula = "=25.95+22.7+14.97+18.72+13.94+13.98+13.44+8.13+4.96+8.29"
$ws.Range("I3").Formula = "=20.53+17.56+16.52+17.61+14.05+13.86+11.2+13.06+6.84+8.05"
$ws.Range("J3").Formula = "=19.56+18.11+16.55+17.13+16.13+15.36+13.8+11+6+8.16"
$ws.Range("K3").Formula = "=20.98+17.64+18.81+20+14.95+13.85+10.74+13.41+8.33+7.96"
$ws.Range("G4").Formula = "=21.34+20.37+16.71+13.71+13.31+11.45+10.59+11.4+6.79+9.14"
$ws.Range("H4").Formula = "=17.39+22.93+8.64+11.77+14.1+9.46+10.04+11.66+6+7"
$ws.Range("I4").Formula = "=21+18.46+14.36+11.79+13.95+10.88+10.22+12.5+6+8.81"
$ws.Range("J4").Formula = "=17.24+20.8+12.98+11.74+15.67+9.61+9.65+12.11+6+8.42"
$ws.Range("K4").Formula = "=16.56+20.66+13.83+12.77+13.97+9.52+9.16+13.13+6+8.76"
$ws.Range("G5").Formula = "=17.41+24.5+10.75+20.6+10.93+8.67+12.21+11.75+6+8.49"
$ws.Range("H5").Formula = "=17.08+23.91+11.5+20.52+13.58+8.75+9.88+10.99+6+8"
$ws.Range("I5").Formula = "=16.91+22.82+11.49+20.89+12.56+8.23+10.74+10.01+6+9.3"
$ws.Range("J5").Formula = "=17.49+27.9+14.11+19.56+11.7+8.45+11.06+11.1+6.68+8.76"
$ws.Range("K5").Formula = "=20.06+25.07+15.39+20.55+12.23+11.29+9.6+10.35+8.11+9.17"
$ws.Range("G6").Formula = "=22.26+16.77+15.19+15.91+13.11+11.05+11.57+13.23+9.04+8.73"
$ws.Range("H6").Formula = "=19.76+16.87+13.47+15.08+12.92+9.75+11.84+8.89+5.14+8.22"
$ws.Range("I6").Formula = "=21.37+14.26+13.12+15.32+13.465+9.88+12.96+13.55+8.41+8.66"
$ws.Range("J6").Formula = "=24.57+15.39+15.35+17.8+13.7+11.47+13.35+13.54+6+8.72"
$ws.Range("K6").Formula = "=21.83+14.85+13.11+16.21+13.91+10.44+11.08+12.11+6+8.43"
$ws.Range("G7").Formula = "=22.36+10.52+20.3+19.97+11.54+17.37+14.15+13.46+6+9.06"
$ws.Range("H7").Formula = "=17+12.01+20.83+19.6+10.41+15.85+13.3+13.74+6+8.5"
$ws.Range("I7").Formula = "=21.04+9.41+19.58+15.4+18.43+15.29+11.52+12.78+5.4+8.88"
$ws.Range("J7").Formula = "=19.34+10.78+23.72+17.88+13.65+16.86+11.77+11.19+6+8.4"
$ws.Range("K7").Formula = "=19.57+10.66+21.75+13.63+21.09+18.56+13.59+12.43+5.7+8.2"
$ws.Range("G8").Formula = "=21.64+17.74+12.66+18.12+11.83+8.07+11.56+10.57+6.01+8.07"
$ws.Range("H8").Formula = "=17.37+15.45+16.07+15.56+12.88+7.9+17.18+12.35+6+7.96"
$ws.Range("I8").Formula = "=18.87+16.1+16.48+16.04+11.66+8.08+12.28+18.73+6+8.28"
$ws.Range("J8").Formula = "=21.04+18.33+15.79+17.57+11.65+8.15+11.33+18.19+5.16+8.4"
$ws.Range("K8").Formula = "=17.3+16.25+16.51+15.54+11.77+8.38+12.59+17.23+7.73+8.02"
$ws.Range("G9").Value = 123.73
$ws.Range("H9").Formula = "=13.92+12.67+10.45+13.86+9.99+8.49+9.79+6.37+8.38"
$ws.Range("I9").Formula = "=13.51+14.02+14.24+15.68+14.02+10.23+9.22+11.02+5.07+8.77"
$ws.Range("J9").Formula = "=17.19+12.82+11.79+13.59+16.29+11.28+10.55+10.27+6+8.52"
$ws.Range("K9").Formula = "=20.02+12.33+12.11+13.53+14.9+10.01+9.77+10.56+6.72+8.16"
$ws.Range("G10").Formula = "=17.95+20.17+12.34+21.53+15.1+14.11+12.91+13.1+6+8.84"
$ws.Range("H10").Formula = "=16.92+16.3+14.16+21.19+17.87+8+11.91+13.1+6.86+8.44"
$ws.Range("I10").Formula = "=18.46+19.33+11.96+20.54+17.35+15.39+12.53+14.32+6+8.3"
$ws.Range("J10").Formula = "=19.1+15.99+11.77+21.66+16.39+15.67+12.5+14.69+6+8.81"
$ws.Range("K10").Formula = "=19.25+14.92+12.95+21.01+16.38+15.75+12.64+12.3+6.75+8.26"
$ws.Range("G11").Formula = "=18.84+15.67+20.71+16.44+11.85+14.56+14.62+15.53+5.7+8.47"
$ws.Range("H11").Formula = "=15.35+13.75+18.82+16.7+11.78+14.04+16.61+12.45+5.98+8.27"
$ws.Range("I11").Formula = "=24.14+14.49+20.8+15.89+10.8+14.93+16.17+12.06+5.48+8.54"
$ws.Range("J11").Formula = "=20.43+13.92+20.68+15.3+10.41+13.95+17.42+12.23+6+7.76"
$ws.Range("K11").Formula = "=19.97+20.23+18.04+14.95+11.86+16.12+17.44+11.62+6+8.09"
$ws.Range("G12").Formula = "=19.37+12.77+14.03+16.79+10.24+7.04+13.89+12.5+5.25+9.14"
$ws.Range("H12").Formula = "=18.89+16.89+14.2+16.7+10.16+8.44+11.09+10.52+5.32+9.17"
$ws.Range("I12").Formula = "=19.49+15.72+13.37+16.99+10.69+8.14+14.54+15.18+8.01+9.33"
$ws.Range("J12").Formula = "=16.96+13.3+12.25+15.85+10.12+7.64+15.31+16.84+9.72+9.21"
$ws.Range("K12").Formula = "=18.95+13+14.18+16.75+10.08+7.64+13.96+13.9+7.85+9.54"
$ws.Range("G13").Formula = "=17.22+14.95+11.82+13.69+15.88+8.91+13.58+14.01+5.86+8.36"
$ws.Range("H13").Formula = "=18.42+12.61+11.62+13.56+15.94+9.1+13.72+9.75+6+8.4"
$ws.Range("I13").Formula = "=17.71+13.97+14.2+14.02+15.48+8.63+13.54+13.05+6+8.37"
$ws.Range("J13").Formula = "=21.63+10.85+13.58+14.6+14.69+8.29+13.26+14.14+6+8.3"
$ws.Range("K13").Formula = "=18.79+12.12+12.83+16.17+16.04+8.07+12.88+14.48+6+8.27"

# Step 4: apply 2-decimal number format to I6 (new style)
$ws.Range("I6").NumberFormat = "0.00"

# Step 5: update selection on Stats sheet
$ws.Range("G16").Select()